# Weekly price-sheet update: insert one new observation row for
# "Femacal de La Calera" / Espinaca ahead of the existing row 358
# (2021-08-17 entry), shifting all subsequent rows down by one and
# extending the used range from A1:R384 to A1:R385.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 358 (rows 358..384 shift to 359..385)
$ws.Rows.Item(358).Insert()

# Fill the new row 358 with the new weekly observation
$ws.Cells.Item(358, 1).Value  = 3
$ws.Cells.Item(358, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(358, 3).Value  = "Coquimbo"
$ws.Cells.Item(358, 4).Value  = 44783
$ws.Cells.Item(358, 5).Value  = 5
$ws.Cells.Item(358, 6).Value  = 100112012
$ws.Cells.Item(358, 7).Value  = "Espinaca"
$ws.Cells.Item(358, 8).Value  = "Sin especificar"
$ws.Cells.Item(358, 9).Value  = "Primera"
$ws.Cells.Item(358, 10).Value = 200
$ws.Cells.Item(358, 11).Value = 3800
$ws.Cells.Item(358, 12).Value = 4000
$ws.Cells.Item(358, 13).Value = 3910
$ws.Cells.Item(358, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(358, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(358, 16).Value = 1303
$ws.Cells.Item(358, 17).Value = 3
$ws.Cells.Item(358, 18).Value = "Hortaliza"
